{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n\n// 1) Remove the entire duplicated second-student block (from the second\n//    \"UNIVERSIDAD CAT\u00d3LICA DEL MAULE\" heading through the page-break that\n//    immediately precedes \"(A COMPLETAR POR PROFESOR GUIA)\"). We locate it\n//    by finding the second occurrence of the heading text and the\n//    paragraph holding \"(A COMPLETAR POR PROFESOR GUIA)\", then delete every\n//    paragraph in between (inclusive on the start side, exclusive of the\n//    \"PROFESOR GUIA\" paragraph itself).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet headingIdx = [];\nlet profesorGuiaIdx = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t === \"UNIVERSIDAD CAT\u00d3LICA DEL MAULE\") {\n    headingIdx.push(i);\n  }\n  if (t === \"(A COMPLETAR POR PROFESOR GUIA)\" && profesorGuiaIdx === -1) {\n    profesorGuiaIdx = i;\n  }\n}\n\nif (headingIdx.length >= 2 && profesorGuiaIdx !== -1) {\n  const startIdx = headingIdx[1]; // second occurrence begins the block to remove\n  const endIdx = profesorGuiaIdx - 1; // paragraph right before \"(A COMPLETAR POR PROFESOR GUIA)\"\n  for (let i = endIdx; i >= startIdx; i--) {\n    paragraphs.items[i].delete();\n  }\n  await context.sync();\n}\n\n// 2) Simple text replacements (run-level text swaps; formatting is kept\n//    because we replace only the matched text range).\nconst replacements = [\n  [\"NOMBRE COMPLETO:Alfonso Bilocopetiuc Parra\", \"NOMBRE COMPLETO:Barbara Alejandra Suarez Sepulveda\"],\n  [\"RUT:18.674.261-3            A\u00d1O INGRESO:2013\", \"RUT:19.123.679-0            A\u00d1O INGRESO:2013\"],\n  [\"EMAIL:  AlfonsoB@gmail.com            TEL\u00c9FONO: 12345678\", \"EMAIL:  BarbaraSuarez@gmail.com            TEL\u00c9FONO: 75632982\"],\n  [\"NOMBRE TESIS/MEMORIA: Tesis 100\", \"NOMBRE TESIS/MEMORIA: Sistema de portal del alumno UCM\"],\n  [\"PROFESOR GU\u00cdA: Marco Toranzo Cespedes\", \"PROFESOR GU\u00cdA: Hugo Araya Carrasco\"],\n  [\"1-Hugo Araya Carrasco\", \"1-Paulo Gonzalez\"],\n  [\"2-Paulo Gonzalez\", \"2-Wladimir Soto\"],\n  [\"3-Ninguno\", \"3-Angelica Urrutia\"],\n  [\"CORREO:JoseT@gmail.com\", \"CORREO:JoseT@utal.com\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, \"Replace\");\n    await context.sync();\n  }\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\n# Helper: collect the Range.Start of every match of $text in the document.\nfunction Find-AllStarts($doc, $text) {\n    $rng = $doc.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Text = $text\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = 0\n    $starts = @()\n    while ($rng.Find.Execute()) {\n        $starts += $rng.Start\n        $rng.Collapse(0)\n    }\n    return $starts\n}\n\n# 1) Remove the entire duplicated second-student block: from the second\n#    \"UNIVERSIDAD CAT\u00d3LICA DEL MAULE\" heading through the page break that\n#    immediately precedes \"(A COMPLETAR POR PROFESOR GUIA)\".\n$headingStarts = Find-AllStarts $d \"UNIVERSIDAD CAT\u00d3LICA DEL MAULE\"\n$guiaStarts = Find-AllStarts $d \"(A COMPLETAR POR PROFESOR GUIA)\"\n\nif ($headingStarts.Count -ge 2 -and $guiaStarts.Count -ge 1) {\n    $deleteStart = $headingStarts[1]\n    $deleteEnd = $guiaStarts[0]\n    $blockRange = $d.Range($deleteStart, $deleteEnd)\n    [void]$blockRange.Delete()\n}\n\n# 2) Simple text replacements (run-level text swaps; formatting is kept\n#    because Find/Replace only rewrites the matched text).\nfunction Replace-Text($doc, $old, $new) {\n    $rng = $doc.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    [void]$rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 1)\n}\n\nReplace-Text $d \"NOMBRE COMPLETO:Alfonso Bilocopetiuc Parra\" \"NOMBRE COMPLETO:Barbara Alejandra Suarez Sepulveda\"\nReplace-Text $d \"RUT:18.674.261-3            A\u00d1O INGRESO:2013\" \"RUT:19.123.679-0            A\u00d1O INGRESO:2013\"\nReplace-Text $d \"EMAIL:  AlfonsoB@gmail.com            TEL\u00c9FONO: 12345678\" \"EMAIL:  BarbaraSuarez@gmail.com            TEL\u00c9FONO: 75632982\"\nReplace-Text $d \"NOMBRE TESIS/MEMORIA: Tesis 100\" \"NOMBRE TESIS/MEMORIA: Sistema de portal del alumno UCM\"\nReplace-Text $d \"PROFESOR GU\u00cdA: Marco Toranzo Cespedes\" \"PROFESOR GU\u00cdA: Hugo Araya Carrasco\"\nReplace-Text $d \"1-Hugo Araya Carrasco\" \"1-Paulo Gonzalez\"\nReplace-Text $d \"2-Paulo Gonzalez\" \"2-Wladimir Soto\"\nReplace-Text $d \"3-Ninguno\" \"3-Angelica Urrutia\"\nReplace-Text $d \"CORREO:JoseT@gmail.com\" \"CORREO:JoseT@utal.com\"\n"}
